# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" column values (column G), replacing the old Strike# derived values.
$kValues = @(4, 7, 6, 3, 3, 3, 8, 4, 6, 6, 5, 5, 1, 9, 4, 5, 5, 7, 7, 5, 6, 2, 1, 8, 6, 4, 7, 4, 10, 3, 4, 5, 3, 6, 3, 4, 2)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
